$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Formula = '="57.427.14"'
$ws.Range('D2').Copy()
$ws.Range('D2').PasteSpecial(-4163)
$ws.Range('E2').Value = '  +2.33%  '

$ws.Range('D3').Formula = '="2.355.46"'
$ws.Range('D3').Copy()
$ws.Range('D3').PasteSpecial(-4163)
$ws.Range('E3').Value = '  +1.77%  '

$ws.Range('E4').Value = '  +0.02%  '

$ws.Range('D5').Formula = '="521.91"'
$ws.Range('D5').Copy()
$ws.Range('D5').PasteSpecial(-4163)
$ws.Range('E5').Value = '  +0.80%  '

$ws.Range('D6').Formula = '="137.18"'
$ws.Range('D6').Copy()
$ws.Range('D6').PasteSpecial(-4163)
$ws.Range('E6').Value = '  +2.97%  '

$ws.Range('D7').Formula = '="0.997"'
$ws.Range('D7').Copy()
$ws.Range('D7').PasteSpecial(-4163)
$ws.Range('E7').Value = '  +0.32%  '

$ws.Range('D8').Formula = '="0.539"'
$ws.Range('D8').Copy()
$ws.Range('D8').PasteSpecial(-4163)
$ws.Range('E8').Value = '  +0.79%  '

$ws.Range('D9').Formula = '="2.374.42"'
$ws.Range('D9').Copy()
$ws.Range('D9').PasteSpecial(-4163)
$ws.Range('E9').Value = '  +1.59%  '

$ws.Range('E10').Value = '  -0.62%  '

$ws.Range('D11').Formula = '="5.46"'
$ws.Range('D11').Copy()
$ws.Range('D11').PasteSpecial(-4163)
$ws.Range('E11').Value = '  +5.72%  '

$ws.Range('E12').Value = '  -1.29%  '

$ws.Range('E13').Value = '  +0.48%  '

$ws.Range('D14').Formula = '="24.24"'
$ws.Range('D14').Copy()
$ws.Range('D14').PasteSpecial(-4163)
$ws.Range('E14').Value = '  +0.89%  '

$ws.Range('D15').Formula = '="2.779.36"'
$ws.Range('D15').Copy()
$ws.Range('D15').PasteSpecial(-4163)
$ws.Range('E15').Value = '  +1.80%  '

$ws.Range('D16').Formula = '="57.444.83"'
$ws.Range('D16').Copy()
$ws.Range('D16').PasteSpecial(-4163)
$ws.Range('E16').Value = '  +2.09%  '

$ws.Range('D17').Formula = '="0.0000136"'
$ws.Range('D17').Copy()
$ws.Range('D17').PasteSpecial(-4163)
$ws.Range('E17').Value = '  +0.47%  '

$ws.Range('D18').Formula = '="2.374.24"'
$ws.Range('D18').Copy()
$ws.Range('D18').PasteSpecial(-4163)
$ws.Range('E18').Value = '  +2.49%  '

$ws.Range('D19').Formula = '="10.64"'
$ws.Range('D19').Copy()
$ws.Range('D19').PasteSpecial(-4163)
$ws.Range('E19').Value = '  +0.89%  '

$ws.Range('D20').Formula = '="329.66"'
$ws.Range('D20').Copy()
$ws.Range('D20').PasteSpecial(-4163)
$ws.Range('E20').Value = '  +2.52%  '

$ws.Range('D21').Formula = '="4.25"'
$ws.Range('D21').Copy()
$ws.Range('D21').PasteSpecial(-4163)
$ws.Range('E21').Value = '  -1.07%  '

$ws.Range('D22').Formula = '="6.74"'
$ws.Range('D22').Copy()
$ws.Range('D22').PasteSpecial(-4163)
$ws.Range('E22').Value = '  +1.07%  '

$ws.Range('D23').Formula = '="0.998"'
$ws.Range('D23').Copy()
$ws.Range('D23').PasteSpecial(-4163)
$ws.Range('E23').Value = '  -0.22%  '

$ws.Range('D24').Formula = '="61.36"'
$ws.Range('D24').Copy()
$ws.Range('D24').PasteSpecial(-4163)
$ws.Range('E24').Value = '  +1.06%  '

$ws.Range('E25').Value = '  +4.43%  '

$ws.Range('D26').Formula = '="0.994"'
$ws.Range('D26').Copy()
$ws.Range('D26').PasteSpecial(-4163)

$ws.Range('D27').Formula = '="8.24"'
$ws.Range('D27').Copy()
$ws.Range('D27').PasteSpecial(-4163)
$ws.Range('E27').Value = '  +7.62%  '

$ws.Range('E28').Value = '  +11.11%  '

$ws.Range('B29').Value = 'Monero'
$ws.Range('C29').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D29').Formula = '="170.76"'
$ws.Range('D29').Copy()
$ws.Range('D29').PasteSpecial(-4163)
$ws.Range('E29').Value = '  -0.17%  '

$ws.Range('B30').Value = 'PEPE'
$ws.Range('C30').Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range('D30').Formula = '="0.0₃0748"'
$ws.Range('D30').Copy()
$ws.Range('D30').PasteSpecial(-4163)
$ws.Range('E30').Value = '  +1.91%  '

$ws.Range('D31').Formula = '="1.71"'
$ws.Range('D31').Copy()
$ws.Range('D31').PasteSpecial(-4163)
$ws.Range('E31').Value = '  -0.74%  '

$ws.Range('D32').Formula = '="6.31"'
$ws.Range('D32').Copy()
$ws.Range('D32').PasteSpecial(-4163)
$ws.Range('E32').Value = '  +0.61%  '

$ws.Range('D33').Formula = '="18.66"'
$ws.Range('D33').Copy()
$ws.Range('D33').PasteSpecial(-4163)
$ws.Range('E33').Value = '  +1.72%  '

$ws.Range('E34').Value = '  +0.03%  '

$ws.Range('E35').Value = '  +2.52%  '

$ws.Range('D36').Formula = '="0.995"'
$ws.Range('D36').Copy()
$ws.Range('D36').PasteSpecial(-4163)
$ws.Range('E36').Value = '  +0.28%  '

$ws.Range('D37').Formula = '="0.929"'
$ws.Range('D37').Copy()
$ws.Range('D37').PasteSpecial(-4163)
$ws.Range('E37').Value = '  +0.29%  '

$ws.Range('D38').Formula = '="4.06"'
$ws.Range('D38').Copy()
$ws.Range('D38').PasteSpecial(-4163)
$ws.Range('E38').Value = '  +1.33%  '

$ws.Range('D39').Formula = '="1.58"'
$ws.Range('D39').Copy()
$ws.Range('D39').PasteSpecial(-4163)
$ws.Range('E39').Value = '  +4.26%  '

$ws.Range('D40').Formula = '="38.54"'
$ws.Range('D40').Copy()
$ws.Range('D40').PasteSpecial(-4163)
$ws.Range('E40').Value = '  +2.95%  '

$ws.Range('D41').Formula = '="151.52"'
$ws.Range('D41').Copy()
$ws.Range('D41').PasteSpecial(-4163)
$ws.Range('E41').Value = '  +8.55%  '

$ws.Range('E42').Value = '  +1.14%  '

$ws.Range('D43').Formula = '="3.65"'
$ws.Range('D43').Copy()
$ws.Range('D43').PasteSpecial(-4163)
$ws.Range('E43').Value = '  +2.09%  '

$ws.Range('D44').Formula = '="5.33"'
$ws.Range('D44').Copy()
$ws.Range('D44').PasteSpecial(-4163)
$ws.Range('E44').Value = '  +3.27%  '

$ws.Range('D45').Formula = '="282.63"'
$ws.Range('D45').Copy()
$ws.Range('D45').PasteSpecial(-4163)
$ws.Range('E45').Value = '  +2.61%  '

$ws.Range('D46').Formula = '="0.0940"'
$ws.Range('D46').Copy()
$ws.Range('D46').PasteSpecial(-4163)
$ws.Range('E46').Value = '  +1.03%  '

$ws.Range('D47').Formula = '="0.0509"'
$ws.Range('D47').Copy()
$ws.Range('D47').PasteSpecial(-4163)
$ws.Range('E47').Value = '  -0.12%  '

$ws.Range('D48').Formula = '="0.566"'
$ws.Range('D48').Copy()
$ws.Range('D48').PasteSpecial(-4163)
$ws.Range('E48').Value = '  +1.55%  '

$ws.Range('E49').Value = '  +2.44%  '

$ws.Range('B50').Value = 'InjectiveProtocol'
$ws.Range('C50').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D50').Formula = '="18.13"'
$ws.Range('D50').Copy()
$ws.Range('D50').PasteSpecial(-4163)
$ws.Range('E50').Value = '  +5.96%  '

$ws.Range('B51').Value = 'Polygon'
$ws.Range('C51').Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range('D51').Formula = '="0.385"'
$ws.Range('D51').Copy()
$ws.Range('D51').PasteSpecial(-4163)
$ws.Range('E51').Value = '  +0.75%  '

$excel.CutCopyMode = 0
